# Quantities.xlsx — shift the date series forward and append two new days.
#
# Rows 2-35: column A is a simple "+2 days" shift (B..J unchanged).
# Rows 36-45: the tail of the series gets reshuffled — two of the existing
# rows change by more than +2 (dates get re-sorted), the two "special"
# (weekend) rows rotate to a new position, and two brand-new rows are
# appended at the bottom (44, 45). Those rows need both the date (A) and,
# where applicable, the dependent quantities (C, G, J) rewritten explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-35: plain "+2 days" shift on column A -------------------------
for ($r = 2; $r -le 35; $r++) {
    $ws.Cells.Item($r, 1).Value = $ws.Cells.Item($r, 1).Value2 + 2
}

# --- Rows 36-45: explicit rewrite of A (and C/G/J where they differ) ------
# row -> (date, C, G, J)
$tailData = @{
    36 = @(45567, 0.00170247, 465.80531254, 485.38834923)
    37 = @(45566, 0.00170247, 465.80531254, 485.38834923)
    38 = @(45564, 0.00170247, 465.80531254, 485.38834923)
    39 = @(45565, 0.00170247, 465.80531254, 485.38834923)
    40 = @(45563, 0.00170247, 465.80531254, 485.38834923)
    41 = @(45558, 0.00004012, 280.99031254, 1941.48834923)
    42 = @(45559, 0.00004012, 280.99031254, 1941.48834923)
    43 = @(45560, 0.00170247, 465.80531254, 485.38834923)
    44 = @(45561, 0.00170247, 465.80531254, 485.38834923)
    45 = @(45562, 0.00170247, 465.80531254, 485.38834923)
}

foreach ($r in 36..45) {
    $vals = $tailData[$r]

    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 7).Value = $vals[2]
    $ws.Cells.Item($r, 10).Value = $vals[3]
}

# Rows 44 and 45 are brand-new — fill in the constant columns (B, D, E, F, H, I)
# that every other row in the table already shares, and mirror row 43's
# date-cell formatting (bold, centered, bordered, yyyy-mm-dd date format) onto
# column A of the two new rows by copying its format (xlPasteFormats).
foreach ($r in 44..45) {
    $ws.Cells.Item($r, 2).Value = 116.4121952
    $ws.Cells.Item($r, 4).Value = 0.00885078
    $ws.Cells.Item($r, 5).Value = 0.06933635
    $ws.Cells.Item($r, 6).Value = 12792.90181321
    $ws.Cells.Item($r, 8).Value = 0.24
    $ws.Cells.Item($r, 9).Value = 1.7904431

    $ws.Range("A43").Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}

Write-Output "done"
